# Insert a new weekly price record at row 7, shifting the existing
# records (previously rows 7-30) down by one row to rows 8-31.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 44761
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 100112012
$ws.Range("G7").Value = "Espinaca"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 35
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("N7").Value = "`$/cuna 10 kilos"
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 1300
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = "Hortaliza"
